# Fixed #418 Empty AQL expressions generate empty lines.
#
# The table in this document has several table cells whose paragraph
# contains nothing but a single "placeholder" run with an empty <w:t/>
# (the leftover of an AQL expression that evaluated to the empty string).
# Those vestigial empty runs should not be present in the generated
# document - the paragraph should simply stay empty (just its <w:pPr>).
#
# For each such cell we can't just clear the run's text (it is already
# empty, so Range.Delete()/Range.Text="" on it is a no-op for the
# underlying empty run) - instead we briefly give the run one character
# of real text (which reuses/keeps the existing run + its <w:pPr>) and
# then delete that character, which collapses/removes the now-truly-empty
# run from the paragraph while leaving the paragraph formatting intact.

function Remove-EmptyCellRun($cell) {
    $para = $cell.Range.Paragraphs.Item(1)
    $rng = $para.Range
    # A single-paragraph cell range's Text always carries 2 "invisible"
    # trailing characters (paragraph mark + cell mark). So Length == 2
    # means the run(s) preceding them hold no real text - i.e. this is
    # one of the "<w:r><w:t/></w:r>" placeholder runs left over from an
    # AQL expression that evaluated to the empty string. Anything longer
    # than 2 is a cell that already has real text - leave it alone.
    if ($rng.Text.Length -ne 2) {
        return
    }
    # Stamp a single placeholder character into the (empty) run...
    $rng.Text = "x"
    # ...then re-acquire the paragraph (content shifted) and delete just
    # that one character, which removes the now-empty run from the XML
    # while leaving the paragraph (and its <w:pPr>) in place.
    $para2 = $cell.Range.Paragraphs.Item(1)
    $rng2 = $para2.Range
    $charRng = $word.ActiveDocument.Range($rng2.Start, $rng2.Start + 1)
    $charRng.Text = ""
}

$d = $word.ActiveDocument

# Walk every table / row / cell in the document and strip out any
# leftover empty-AQL-expression run we find, wherever it is. This does
# not hard-code row/column positions, so it applies to every table cell
# whose paragraph is really just an empty placeholder run.
foreach ($table in $d.Tables) {
    foreach ($row in $table.Rows) {
        foreach ($cell in $row.Cells) {
            Remove-EmptyCellRun $cell
        }
    }
}
